{"js": "// Load all paragraphs in the document body so we can locate the\n// whitespace-only paragraphs near the end of the second code block\n// and the final closing-brace paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the two consecutive paragraphs whose entire text is exactly\n// eight spaces (\"        \") that sit right before the very last\n// \"    }\" paragraph, and delete them.\nlet targetIndex = -1;\nfor (let i = 0; i < items.length - 1; i++) {\n  if (items[i].text === \"        \" && items[i + 1].text === \"        \") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  items[targetIndex].delete();\n  items[targetIndex + 1].delete();\n  await context.sync();\n}\n\n// Re-fetch paragraphs after the deletion so we operate on a fresh,\n// up-to-date collection, then append the new content after the very\n// last paragraph in the document (the final \"    }\").\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nconst refreshedItems = refreshedParagraphs.items;\nlet anchor = refreshedItems[refreshedItems.length - 1];\n\nanchor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\nanchor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\nanchor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\nanchor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\nanchor = anchor.insertParagraph(\"Hello my name is het shah.\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two consecutive whitespace-only paragraphs (each holding\n# exactly eight spaces) that sit right before the final \"    }\"\n# paragraph of the second code listing, and remove them.\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -lt $count; $i++) {\n    $t1 = $d.Paragraphs.Item($i).Range.Text\n    $t2 = $d.Paragraphs.Item($i + 1).Range.Text\n    if ($t1 -eq \"        `r\" -and $t2 -eq \"        `r\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ne -1) {\n    # Deleting the paragraph at $targetIndex twice removes both,\n    # since the next paragraph shifts into $targetIndex after each delete.\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n}\n\n# Append five new paragraphs after the very last paragraph in the\n# document (the closing \"    }\") : four empty ones followed by a\n# paragraph of text.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n$lastPara.Range.InsertParagraphAfter()\n$lastPara.Range.InsertParagraphAfter()\n$lastPara.Range.InsertParagraphAfter()\n$lastPara.Range.InsertParagraphAfter()\n\n$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$finalPara.Range.Text = \"Hello my name is het shah.\"\n"}
